$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Quantité Ingrédients"), pushing the
# existing C..G columns to D..H. This becomes the new "Quantite" (numeric)
# column.
$ws.Range("C1").EntireColumn.Insert()

# Header for the new column
$ws.Range("C1").Value = "Quantite"

# Numeric quantities (previously encoded only as text like "1u", "4u", "100g"
# inside what is now column D)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 100
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 600

# Give the new column the same width as column B
$ws.Range("C1").ColumnWidth = 19

# Move the active selection to C8, matching the saved view state
$ws.Range("C8").Select() | Out-Null
